$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (shifts old rows 13-23 down to 14-24) ---
$ws.Rows("13:13").Insert()

# The freshly inserted row's cells don't reliably inherit the per-column
# style (B -> style 2, C -> style 3), so copy formatting from an existing
# correctly-styled row (row 9: A=1, B=2, C=3) onto the new row 13 first.
$ws.Range("A9:C9").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

# Row 13 has no label in column A in the final layout - clear it so no
# stray empty <c> element remains there.
$ws.Range("A13").Clear()

# --- Fill in new content (B & C mirror each other throughout this sheet) ---

# Row 10 (Objetivos:) now gets its answer text (previously blank here -
# the "Daisy" text that used to sit here moves down into the new row 13).
$objetivosText = "Proporcionar aos estudantes dos cursos de Engenharia da EEL-USP a reflexão crítica sobre o pensamento filosófico e a evolução das ciências no mundo contemporâneo."
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# Row 13 (new row, under "Docentes responsáveis:") - Daisy Rafaela text.
$docenteText = "6376612 - Daisy Rafaela da Silva"
$ws.Range("B13").Value = $docenteText
$ws.Range("C13").Value = $docenteText

# Row 14 (Programa resumido:) gets its answer text (was "Semestral").
$resumidoText = "Filosofia e epistemologia do conhecimento.  As ciências e a metodologia do conhecimento científico.  Modernidade e pós-modernidade. Questões éticas e sócio-ambiental."
$ws.Range("B14").Value = $resumidoText
$ws.Range("C14").Value = $resumidoText

# Row 16 (Programa:) gets its answer text (previously held the
# 01/01/2012 date by mistake).
$programaText = "Unidade primeira:  Filosofia e Ciência`n 1 - A filosofia`n2 -  Formas de conceber o conhecimento.  `n 2  Epistemologia do conhecimento científico`n 3 - A metodologia científica `nUnidade Segunda: Reflexões sobre a ciência `n1  Ciência:  Modernidade e pós-Modernidade`n2 - Ciência e ética`n3 - Ciência e a questão sócio-ambiental `n4 - Ciência, tecnologia e sustentabilidade"
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# Row 19 (Método:) gets its answer text (previously held the Daisy text
# by mistake).
$metodoText = "A média semestral e final dos alunos será composta por: Prova Semestral  (PS) e outros instrumentos (T) empregados na avaliação do aluno, valorizando a sua participação e colaboração nos trabalhos e atividades desenvolvidas individualmente e no Projeto de curso em equipe."
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20 (Critério:) gets its answer text (previously held the "Método"
# answer text by mistake).
$criterioText = "(PS+T) / 2"
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21 (Norma de recuperação:) gets its answer text (previously held
# the "Critério" answer text by mistake).
$normaText = "Trabalho escrito, com questionamento, envolvendo o conteúdo do programa `n-   prova escrita"
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# Row 22 (Bibliografia:) gets its answer text (previously held the
# "Norma de recuperação" answer text by mistake).
$bibliografiaText = "1 - CAPRA, F. O Ponto de Mutação. São Paulo: Cultrix, 1986. `n2 - CAPRA, F. A Teia da Vida.  São Paulo: Cultrix, 2003.`n3- CASTELLS,Manuel. O Poder da  Identidade. A Era da Informação: Economia, Sociedade Cultura. Vol. 2. Rio de Janeiro: Paz e      Terra, 1999. `n4 -CERVO, Amado L.; BERVIAN, Pedro. A Metodologia Científica. São Paulo: Prentice      Hall, 2002. `n5 - POPPER, Karl. R. A Lógica da Pesquisa Científica. São Paulo: Cutrix, 1959`n6 - RAMPAZZO, Lino. Metodologia Científica. São Paulo: Edições Loyola, 3ª ed. , 2005.`n7 - SANTOS, Boaventura de.  Um Discurso sobre as Ciências.  Porto, Portugal: Aforntamentos, 1997.`n8 - STEGMULLER,  W. A Filosofia Contemporânea. Vol. I e II, 1977."
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText
